# Add a new "2022-Q1" fund-holding detail sheet (positioned right after
# "2021-Q4" and right before "总计"), and refresh the "总计" (summary)
# sheet with a new leading row for 2022-Q1 - the previously existing
# rows shift down by one.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# The old summary sheet is rebuilt from scratch right after the new
# detail sheet, so sheet/rel numbering comes out sequential again
# (2022-Q1 -> sheetId 4, 总计 -> sheetId 5), matching a freshly
# re-saved workbook.
$oldTotal.Delete()

# ---------------------------------------------------------------------
# 1) New detail sheet "2022-Q1" (same layout as the other quarterly
#    detail sheets, e.g. 2021-Q3 / 2021-Q4).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Copy the header-row formatting (bold font + border, style used by the
# other detail sheets) from an existing detail sheet so the new sheet
# matches the same look without re-creating styles by hand.
$q4.Range("B1:H1").Copy()
$q1.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Formula = "=""002076"""
$q1.Range("C2").Value = "浙商中证500指数增强A"
$q1.Range("D2").Formula = "=""14.53"""
$q1.Range("E2").Formula = "=""93.68"""
$q1.Range("F2").Formula = "=""1.42"""
$q1.Range("G2").Formula = "=""0.2063"""
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").Formula = "=""007386"""
$q1.Range("C3").Value = "浙商中证500指数增强C"
$q1.Range("D3").Formula = "=""3.38"""
$q1.Range("E3").Formula = "=""93.68"""
$q1.Range("F3").Formula = "=""1.42"""
$q1.Range("G3").Formula = "=""0.0480"""
$q1.Range("H3").Value = 10

# Freeze the helper formulas down into plain text values (no residual
# "Text" number-format gets left behind on the cells this way).
$q1.Range("B2:G3").Copy()
$q1.Range("B2").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2) Rebuild the "总计" summary sheet with the new 2022-Q1 row on top
#    and the previous rows shifted down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$q4.Range("B1:D1").Copy()
$total.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.25

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.45

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.09

# Restore the originally active tab (2021-Q2) - adding/renaming sheets
# above shifts the "active sheet" around.
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name
}
